$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 25: only B changes
$ws.Range("B25").Value = 96720

# Row 26
$ws.Range("A26").Value = 112183134
$ws.Range("B26").Value = 96720
$ws.Range("Q26").Value = 572354
$ws.Range("R26").Value = 6714968

# Row 27
$ws.Range("A27").Value = 112183150
$ws.Range("B27").Value = 96720
$ws.Range("Q27").Value = 572358
$ws.Range("R27").Value = 6714972

# Row 28
$ws.Range("A28").Value = 112183149
$ws.Range("B28").Value = 96720
$ws.Range("Q28").Value = 572345
$ws.Range("R28").Value = 6714965

# Row 29
$ws.Range("A29").Value = 112183151
$ws.Range("B29").Value = 96720
$ws.Range("Q29").Value = 572361
$ws.Range("R29").Value = 6714978

# Row 30
$ws.Range("A30").Value = 112183140
$ws.Range("B30").Value = 96720
$ws.Range("Q30").Value = 572350
$ws.Range("R30").Value = 6714962

# Row 31
$ws.Range("A31").Value = 112183146
$ws.Range("B31").Value = 96720
$ws.Range("Q31").Value = 572346
$ws.Range("R31").Value = 6714917

# Row 32
$ws.Range("A32").Value = 112183147
$ws.Range("B32").Value = 96720
$ws.Range("Q32").Value = 572351
$ws.Range("R32").Value = 6714915

# Row 33
$ws.Range("A33").Value = 112183145
$ws.Range("B33").Value = 96720
$ws.Range("Q33").Value = 572351
$ws.Range("R33").Value = 6714907

# Row 34
$ws.Range("A34").Value = 112183137
$ws.Range("B34").Value = 96720
$ws.Range("Q34").Value = 572354
$ws.Range("R34").Value = 6714961

# Row 35
$ws.Range("A35").Value = 112183143
$ws.Range("B35").Value = 96720
$ws.Range("Q35").Value = 572359
$ws.Range("R35").Value = 6714905

# Row 36
$ws.Range("A36").Value = 112183148
$ws.Range("B36").Value = 96720
$ws.Range("Q36").Value = 572357
$ws.Range("R36").Value = 6714903
